# Auto update Excel log: append newly-logged sensor readings (2026-01-28, 15:32-15:33)
# to the PIR, Humidity and Temperature sheets of the SeniorConnect master log.
$wb = $excel.ActiveWorkbook

# ---- PIR: append rows 54-66 ----
$ws = $wb.Worksheets.Item('PIR')
$ws.Range('A54').NumberFormat = "@"
$ws.Range('A54').Value = '2026-01-28'
$ws.Range('B54').Value = '15:32:48'
$ws.Range('C54').Value = '15:00'
$ws.Range('D54').Value = 'Bathroom'
$ws.Range('E54').Value = 'No Motion'
$ws.Range('F54').Value = 'Inactive'

$ws.Range('A55').NumberFormat = "@"
$ws.Range('A55').Value = '2026-01-28'
$ws.Range('B55').Value = '15:32:49'
$ws.Range('C55').Value = '15:00'
$ws.Range('D55').Value = 'Bathroom'
$ws.Range('E55').Value = 'No Motion'
$ws.Range('F55').Value = 'Inactive'

$ws.Range('A56').NumberFormat = "@"
$ws.Range('A56').Value = '2026-01-28'
$ws.Range('B56').Value = '15:32:54'
$ws.Range('C56').Value = '15:00'
$ws.Range('D56').Value = 'Bathroom'
$ws.Range('E56').Value = 'No Motion'
$ws.Range('F56').Value = 'Inactive'

$ws.Range('A57').NumberFormat = "@"
$ws.Range('A57').Value = '2026-01-28'
$ws.Range('B57').Value = '15:32:59'
$ws.Range('C57').Value = '15:00'
$ws.Range('D57').Value = 'Bathroom'
$ws.Range('E57').Value = 'No Motion'
$ws.Range('F57').Value = 'Inactive'

$ws.Range('A58').NumberFormat = "@"
$ws.Range('A58').Value = '2026-01-28'
$ws.Range('B58').Value = '15:33:04'
$ws.Range('C58').Value = '15:00'
$ws.Range('D58').Value = 'Bathroom'
$ws.Range('E58').Value = 'No Motion'
$ws.Range('F58').Value = 'Inactive'

$ws.Range('A59').NumberFormat = "@"
$ws.Range('A59').Value = '2026-01-28'
$ws.Range('B59').Value = '15:33:09'
$ws.Range('C59').Value = '15:00'
$ws.Range('D59').Value = 'Bathroom'
$ws.Range('E59').Value = 'No Motion'
$ws.Range('F59').Value = 'Inactive'

$ws.Range('A60').NumberFormat = "@"
$ws.Range('A60').Value = '2026-01-28'
$ws.Range('B60').Value = '15:33:14'
$ws.Range('C60').Value = '15:00'
$ws.Range('D60').Value = 'Bathroom'
$ws.Range('E60').Value = 'No Motion'
$ws.Range('F60').Value = 'Inactive'

$ws.Range('A61').NumberFormat = "@"
$ws.Range('A61').Value = '2026-01-28'
$ws.Range('B61').Value = '15:33:19'
$ws.Range('C61').Value = '15:00'
$ws.Range('D61').Value = 'Bathroom'
$ws.Range('E61').Value = 'No Motion'
$ws.Range('F61').Value = 'Inactive'

$ws.Range('A62').NumberFormat = "@"
$ws.Range('A62').Value = '2026-01-28'
$ws.Range('B62').Value = '15:33:24'
$ws.Range('C62').Value = '15:00'
$ws.Range('D62').Value = 'Bathroom'
$ws.Range('E62').Value = 'No Motion'
$ws.Range('F62').Value = 'Inactive'

$ws.Range('A63').NumberFormat = "@"
$ws.Range('A63').Value = '2026-01-28'
$ws.Range('B63').Value = '15:33:29'
$ws.Range('C63').Value = '15:00'
$ws.Range('D63').Value = 'Bathroom'
$ws.Range('E63').Value = 'No Motion'
$ws.Range('F63').Value = 'Inactive'

$ws.Range('A64').NumberFormat = "@"
$ws.Range('A64').Value = '2026-01-28'
$ws.Range('B64').Value = '15:33:34'
$ws.Range('C64').Value = '15:00'
$ws.Range('D64').Value = 'Bathroom'
$ws.Range('E64').Value = 'No Motion'
$ws.Range('F64').Value = 'Inactive'

$ws.Range('A65').NumberFormat = "@"
$ws.Range('A65').Value = '2026-01-28'
$ws.Range('B65').Value = '15:33:39'
$ws.Range('C65').Value = '15:00'
$ws.Range('D65').Value = 'Bathroom'
$ws.Range('E65').Value = 'No Motion'
$ws.Range('F65').Value = 'Inactive'

$ws.Range('A66').NumberFormat = "@"
$ws.Range('A66').Value = '2026-01-28'
$ws.Range('B66').Value = '15:33:44'
$ws.Range('C66').Value = '15:00'
$ws.Range('D66').Value = 'Bathroom'
$ws.Range('E66').Value = 'No Motion'
$ws.Range('F66').Value = 'Inactive'

# ---- Humidity: append rows 58-67 ----
$ws = $wb.Worksheets.Item('Humidity')
$ws.Range('A58').NumberFormat = "@"
$ws.Range('E58').NumberFormat = "@"
$ws.Range('A58').Value = '2026-01-28'
$ws.Range('B58').Value = '15:32:54'
$ws.Range('C58').Value = '15:00'
$ws.Range('D58').Value = 'Bathroom'
$ws.Range('E58').Value = '87.3%'
$ws.Range('F58').Value = 'Active'

$ws.Range('A59').NumberFormat = "@"
$ws.Range('E59').NumberFormat = "@"
$ws.Range('A59').Value = '2026-01-28'
$ws.Range('B59').Value = '15:32:58'
$ws.Range('C59').Value = '15:00'
$ws.Range('D59').Value = 'Bathroom'
$ws.Range('E59').Value = '88.3%'
$ws.Range('F59').Value = 'Active'

$ws.Range('A60').NumberFormat = "@"
$ws.Range('E60').NumberFormat = "@"
$ws.Range('A60').Value = '2026-01-28'
$ws.Range('B60').Value = '15:33:02'
$ws.Range('C60').Value = '15:00'
$ws.Range('D60').Value = 'Bathroom'
$ws.Range('E60').Value = '87.3%'
$ws.Range('F60').Value = 'Active'

$ws.Range('A61').NumberFormat = "@"
$ws.Range('E61').NumberFormat = "@"
$ws.Range('A61').Value = '2026-01-28'
$ws.Range('B61').Value = '15:33:06'
$ws.Range('C61').Value = '15:00'
$ws.Range('D61').Value = 'Bathroom'
$ws.Range('E61').Value = '88.2%'
$ws.Range('F61').Value = 'Active'

$ws.Range('A62').NumberFormat = "@"
$ws.Range('E62').NumberFormat = "@"
$ws.Range('A62').Value = '2026-01-28'
$ws.Range('B62').Value = '15:33:14'
$ws.Range('C62').Value = '15:00'
$ws.Range('D62').Value = 'Bathroom'
$ws.Range('E62').Value = '88.2%'
$ws.Range('F62').Value = 'Active'

$ws.Range('A63').NumberFormat = "@"
$ws.Range('E63').NumberFormat = "@"
$ws.Range('A63').Value = '2026-01-28'
$ws.Range('B63').Value = '15:33:18'
$ws.Range('C63').Value = '15:00'
$ws.Range('D63').Value = 'Bathroom'
$ws.Range('E63').Value = '88.2%'
$ws.Range('F63').Value = 'Active'

$ws.Range('A64').NumberFormat = "@"
$ws.Range('E64').NumberFormat = "@"
$ws.Range('A64').Value = '2026-01-28'
$ws.Range('B64').Value = '15:33:22'
$ws.Range('C64').Value = '15:00'
$ws.Range('D64').Value = 'Bathroom'
$ws.Range('E64').Value = '87.4%'
$ws.Range('F64').Value = 'Active'

$ws.Range('A65').NumberFormat = "@"
$ws.Range('E65').NumberFormat = "@"
$ws.Range('A65').Value = '2026-01-28'
$ws.Range('B65').Value = '15:33:26'
$ws.Range('C65').Value = '15:00'
$ws.Range('D65').Value = 'Bathroom'
$ws.Range('E65').Value = '88.3%'
$ws.Range('F65').Value = 'Active'

$ws.Range('A66').NumberFormat = "@"
$ws.Range('E66').NumberFormat = "@"
$ws.Range('A66').Value = '2026-01-28'
$ws.Range('B66').Value = '15:33:38'
$ws.Range('C66').Value = '15:00'
$ws.Range('D66').Value = 'Bathroom'
$ws.Range('E66').Value = '86.9%'
$ws.Range('F66').Value = 'Active'

$ws.Range('A67').NumberFormat = "@"
$ws.Range('E67').NumberFormat = "@"
$ws.Range('A67').Value = '2026-01-28'
$ws.Range('B67').Value = '15:33:42'
$ws.Range('C67').Value = '15:00'
$ws.Range('D67').Value = 'Bathroom'
$ws.Range('E67').Value = '87.5%'
$ws.Range('F67').Value = 'Active'

# ---- Temperature: append rows 58-67 ----
$ws = $wb.Worksheets.Item('Temperature')
$ws.Range('A58').NumberFormat = "@"
$ws.Range('A58').Value = '2026-01-28'
$ws.Range('B58').Value = '15:32:54'
$ws.Range('C58').Value = '15:00'
$ws.Range('D58').Value = 'Bathroom'
$ws.Range('E58').Value = '22.9C'
$ws.Range('F58').Value = 'Active'

$ws.Range('A59').NumberFormat = "@"
$ws.Range('A59').Value = '2026-01-28'
$ws.Range('B59').Value = '15:32:58'
$ws.Range('C59').Value = '15:00'
$ws.Range('D59').Value = 'Bathroom'
$ws.Range('E59').Value = '22.9C'
$ws.Range('F59').Value = 'Active'

$ws.Range('A60').NumberFormat = "@"
$ws.Range('A60').Value = '2026-01-28'
$ws.Range('B60').Value = '15:33:02'
$ws.Range('C60').Value = '15:00'
$ws.Range('D60').Value = 'Bathroom'
$ws.Range('E60').Value = '22.9C'
$ws.Range('F60').Value = 'Active'

$ws.Range('A61').NumberFormat = "@"
$ws.Range('A61').Value = '2026-01-28'
$ws.Range('B61').Value = '15:33:06'
$ws.Range('C61').Value = '15:00'
$ws.Range('D61').Value = 'Bathroom'
$ws.Range('E61').Value = '22.9C'
$ws.Range('F61').Value = 'Active'

$ws.Range('A62').NumberFormat = "@"
$ws.Range('A62').Value = '2026-01-28'
$ws.Range('B62').Value = '15:33:14'
$ws.Range('C62').Value = '15:00'
$ws.Range('D62').Value = 'Bathroom'
$ws.Range('E62').Value = '22.9C'
$ws.Range('F62').Value = 'Active'

$ws.Range('A63').NumberFormat = "@"
$ws.Range('A63').Value = '2026-01-28'
$ws.Range('B63').Value = '15:33:18'
$ws.Range('C63').Value = '15:00'
$ws.Range('D63').Value = 'Bathroom'
$ws.Range('E63').Value = '22.9C'
$ws.Range('F63').Value = 'Active'

$ws.Range('A64').NumberFormat = "@"
$ws.Range('A64').Value = '2026-01-28'
$ws.Range('B64').Value = '15:33:22'
$ws.Range('C64').Value = '15:00'
$ws.Range('D64').Value = 'Bathroom'
$ws.Range('E64').Value = '22.9C'
$ws.Range('F64').Value = 'Active'

$ws.Range('A65').NumberFormat = "@"
$ws.Range('A65').Value = '2026-01-28'
$ws.Range('B65').Value = '15:33:26'
$ws.Range('C65').Value = '15:00'
$ws.Range('D65').Value = 'Bathroom'
$ws.Range('E65').Value = '22.9C'
$ws.Range('F65').Value = 'Active'

$ws.Range('A66').NumberFormat = "@"
$ws.Range('A66').Value = '2026-01-28'
$ws.Range('B66').Value = '15:33:38'
$ws.Range('C66').Value = '15:00'
$ws.Range('D66').Value = 'Bathroom'
$ws.Range('E66').Value = '22.9C'
$ws.Range('F66').Value = 'Active'

$ws.Range('A67').NumberFormat = "@"
$ws.Range('A67').Value = '2026-01-28'
$ws.Range('B67').Value = '15:33:42'
$ws.Range('C67').Value = '15:00'
$ws.Range('D67').Value = 'Bathroom'
$ws.Range('E67').Value = '22.9C'
$ws.Range('F67').Value = 'Active'
